# Add the REMIND-MAgPIE coupled model registration as a new row (row 6)
# to the IPCC AR6 Model Registrations overview sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a "dd.mm.yyyy"-style date recorded as plain text (as the
# existing rows do), so force text formatting before/after the write to
# avoid Excel auto-converting the string into a date serial number while
# still keeping the cell on the sheet's default (unstyled) format.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "10.06.2020"
$ws.Range("A6").NumberFormat = "General"

$ws.Range("B6").Value = "REMIND-MAgPIE 2.1-4.2"
$ws.Range("C6").Value = "IPCC_AR6_model_registration_REMIND-MAgPIE_2.1-4.2.xlsx"
$ws.Range("D6").Value = "Björn Sörgel"
$ws.Range("E6").Value = "REMIND – MAgPIE Coupled Version."

# Move the active selection down below the newly added row, matching the
# author's saved cursor position after entering the new data.
$ws.Range("A7").Select() | Out-Null
